# Weekly data refresh: insert 4 new price rows for the most recent report
# date (2021-09-29) ahead of the existing historical rows, which shift down
# by four rows (234-301 -> 238-305).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 234..301 down by 4 rows to make room for the new
# weekly entries (mirrors Excel's Rows.Insert, which also grows the sheet's
# used-range / dimension automatically).
$ws.Rows("234:237").Insert()

$newRows = @(
    @{ Row=234; D=44468; K="Lane Late";  L="Primera"; M=20; N=120000; O=130000; P=125000; S=312 },
    @{ Row=235; D=44468; K="Lane Late";  L="Segunda"; M=20; N=100000; O=110000; P=105000; S=262 },
    @{ Row=236; D=44468; K="Navel Late"; L="Primera"; M=28; N=120000; O=130000; P=126429; S=316 },
    @{ Row=237; D=44468; K="Navel Late"; L="Segunda"; M=20; N=100000; O=110000; P=105000; S=262 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 2
    $ws.Cells.Item($row, 2).Value  = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = 4
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100102
    $ws.Cells.Item($row, 8).Value  = "Cítricos"
    $ws.Cells.Item($row, 9).Value  = 100102005
    $ws.Cells.Item($row, 10).Value = "Naranja"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/bins (400 kilos)"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 400
}
